# Registration Page and Functions
# Populates the "Registration" sheet with a full set of registration test
# data (FirstName/LastName/Phone/Email/Address/... columns), wires up the
# mailto: hyperlink for the Email cell, and makes "Registration" the
# active/selected sheet (it was "CreateAccount" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

# ---------------------------------------------------------------------
# Header row (row 1) -- A1/B1 ("Flag"/"TestCaseName") already exist, add
# the rest of the headers C1:M1.
# ---------------------------------------------------------------------
$ws.Range("C1").Value = "FirstName"
$ws.Range("D1").Value = "LastName"
$ws.Range("E1").Value = "Phone"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Address"
$ws.Range("H1").Value = "City"
$ws.Range("I1").Value = "State"
$ws.Range("J1").Value = "PostalCode"
$ws.Range("K1").Value = "Country"
$ws.Range("L1").Value = "UserName"
$ws.Range("M1").Value = "Password"

# Give the new header cells the same "header" look (fill + border) as the
# existing A1/B1 cells by copying their format across.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Data row (row 2) -- A2/B2 ("Yes"/"TC01") already exist, add the rest of
# the registration record C2:M2.
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "Raj"
$ws.Range("D2").Value = "Sharma"
# Phone/PostalCode look like numbers -- prefix with an apostrophe so they
# are stored as text (quotePrefix), same as Excel does when you type them.
$ws.Range("E2").Value = "'7761237786"
$ws.Range("G2").Value = "1506, VeteranPkway"
$ws.Range("H2").Value = "Macon"
$ws.Range("I2").Value = "Georgia"
$ws.Range("J2").Value = "'31670"
$ws.Range("K2").Value = "UNITED STATES"
$ws.Range("L2").Value = "raj_sharma"
$ws.Range("M2").Value = "Welcome@123"

# Email cell: set the display text then turn it into a real mailto: link.
$ws.Range("F2").Value = "raj_sharma@email.com"
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:raj_sharma@email.com") | Out-Null

# ---------------------------------------------------------------------
# Column widths -- autofit the new columns similar to the existing ones.
# ---------------------------------------------------------------------
for ($col = 1; $col -le 13; $col++) {
    if ($col -eq 2) { continue }  # column B width is already set/bestFit
    $ws.Columns.Item($col).AutoFit() | Out-Null
}

# ---------------------------------------------------------------------
# Make "Registration" the active sheet/tab (it replaces "CreateAccount").
# ---------------------------------------------------------------------
$ws.Range("L3").Select() | Out-Null
$ws.Activate()

Write-Output "Registration sheet populated"
